$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# --- Title shape: reposition/resize, enable shrink-to-fit, update text ---
$title = $s.Shapes.Item("Title 2")

$title.Left = 50.80952755905512
$title.Top = 223.56001
$title.Width = 487.43999
$title.Height = 92.88001

$title.TextFrame.TextRange.Text = "Product Goal + DoD"
$title.TextFrame.AutoSize = 2

# --- Remove the Content/Date/Footer placeholders (kept Slide Number + Picture) ---
$s.Shapes.Item("Content Placeholder 3").Cut()
$s.Shapes.Item("Date Placeholder 8").Cut()
$s.Shapes.Item("Footer Placeholder 9").Cut()
